# Update "想去人数" (want-to-go count) figures for a couple of events.
# Sheet "展览" (Exhibition) and Sheet "全部类型" (All types) both list the
# same events, so both need to be updated in lockstep.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 321   # 南宁·第五届小蜜蜂动漫嘉年华: 320 -> 321
$ws1.Range("F4").Value = 1294  # 南宁·草莓动漫节: 1292 -> 1294
$ws1.Range("F5").Value = 640   # 南宁·第一届ANE·DACG动漫嘉年华: 639 -> 640

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 321   # 南宁·第五届小蜜蜂动漫嘉年华: 320 -> 321
$ws4.Range("F4").Value = 1294  # 南宁·草莓动漫节: 1292 -> 1294
$ws4.Range("F6").Value = 640   # 南宁·第一届ANE·DACG动漫嘉年华: 639 -> 640
